$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare style templates in scratch cells (style 3 = general centered cell, style 4 = price cell) ---
# Row 29 (pre-edit) still has untouched original styling (s=3 for A/B/C/E/F, s=4 for D).
$ws.Cells.Item(29,5).Copy() | Out-Null
$ws.Cells.Item(1,10).PasteSpecial(-4122) | Out-Null   # style-3 template (col J=10)
$ws.Cells.Item(29,4).Copy() | Out-Null
$ws.Cells.Item(1,11).PasteSpecial(-4122) | Out-Null   # style-4 template (col K=11)
$excel.CutCopyMode = $false

# --- Remove existing hyperlinks (they will be re-created at their new, shifted locations) ---
$ws.Hyperlinks.Delete()

# --- Write Sl.no. (A), Description (B), Grade (C), Basic Price (D), Circular Date (E), Circular Link (F) ---
$desc = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$grade = "P1020"

$ws.Cells.Item(2,1).Value2 = 29
$ws.Cells.Item(2,2).Value2 = $desc
$ws.Cells.Item(2,3).Value2 = $grade
$ws.Cells.Item(2,4).Value2 = 263
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value2 = "02.09.2025"
$ws.Cells.Item(2,6).Value2 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-02-september-2025.pdf"

$ws.Cells.Item(3,1).Value2 = 28
$ws.Cells.Item(3,2).Value2 = $desc
$ws.Cells.Item(3,3).Value2 = $grade
$ws.Cells.Item(3,4).Value2 = 261.25
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value2 = "27.08.2025"
$ws.Cells.Item(3,6).Value2 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-27-august-2025.pdf"

$ws.Cells.Item(4,1).Value2 = 27
$ws.Cells.Item(4,2).Value2 = $desc
$ws.Cells.Item(4,3).Value2 = $grade
$ws.Cells.Item(4,4).Value2 = 258.25
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value2 = "26.08.2025"
$ws.Cells.Item(4,6).Value2 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-26-august-2025.pdf"

$ws.Cells.Item(5,1).Value2 = 26
$ws.Cells.Item(5,2).Value2 = $desc
$ws.Cells.Item(5,3).Value2 = $grade
$ws.Cells.Item(5,4).Value2 = 265
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value2 = "23.08.2025"
$ws.Cells.Item(5,6).Value2 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-august-2025.pdf"

$ws.Cells.Item(6,1).Value2 = 25
$ws.Cells.Item(6,2).Value2 = $desc
$ws.Cells.Item(6,3).Value2 = $grade
$ws.Cells.Item(6,4).Value2 = 262.5
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value2 = "20.08.2025"
$ws.Cells.Item(6,6).Value2 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-august-2025.pdf"

$ws.Cells.Item(7,1).Value2 = 24
$ws.Cells.Item(7,2).Value2 = $desc
$ws.Cells.Item(7,3).Value2 = $grade
$ws.Cells.Item(7,4).Value2 = 264.75
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value2 = "19.08.2025"
$ws.Cells.Item(7,6).Value2 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-19-august-2025.pdf"

$ws.Cells.Item(8,1).Value2 = 23
$ws.Cells.Item(8,2).Value2 = $desc
$ws.Cells.Item(8,3).Value2 = $grade
$ws.Cells.Item(8,4).Value2 = 269.25
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value2 = "14.08.2025"
$ws.Cells.Item(8,6).Value2 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-14-august-2025.pdf"

$ws.Cells.Item(9,1).Value2 = 22
$ws.Cells.Item(9,2).Value2 = $desc
$ws.Cells.Item(9,3).Value2 = $grade
$ws.Cells.Item(9,4).Value2 = 267.25
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value2 = "13.08.2025"
$ws.Cells.Item(9,6).Value2 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-august-2025.pdf"

$ws.Cells.Item(10,1).Value2 = 21
$ws.Cells.Item(10,2).Value2 = $desc
$ws.Cells.Item(10,3).Value2 = $grade
$ws.Cells.Item(10,4).Value2 = 265.25
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value2 = "12.08.2025"
$ws.Cells.Item(10,6).Value2 = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf"

$ws.Cells.Item(11,1).Value2 = 20
$ws.Cells.Item(11,2).Value2 = $desc
$ws.Cells.Item(11,3).Value2 = $grade
$ws.Cells.Item(11,4).Value2 = 268.5
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value2 = "08.08.2025"
$ws.Cells.Item(11,6).Value2 = ""

$ws.Cells.Item(12,1).Value2 = 19
$ws.Cells.Item(12,2).Value2 = $desc
$ws.Cells.Item(12,3).Value2 = $grade
$ws.Cells.Item(12,4).Value2 = 265.75
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value2 = "07.08.2025"
$ws.Cells.Item(12,6).Value2 = ""

$ws.Cells.Item(13,1).Value2 = 18
$ws.Cells.Item(13,2).Value2 = $desc
$ws.Cells.Item(13,3).Value2 = $grade
$ws.Cells.Item(13,4).Value2 = 263.75
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value2 = "05.08.2025"
$ws.Cells.Item(13,6).Value2 = ""

$ws.Cells.Item(14,1).Value2 = 17
$ws.Cells.Item(14,2).Value2 = $desc
$ws.Cells.Item(14,3).Value2 = $grade
$ws.Cells.Item(14,4).Value2 = 260.5
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value2 = "02.08.2025"
$ws.Cells.Item(14,6).Value2 = ""

$ws.Cells.Item(15,1).Value2 = 16
$ws.Cells.Item(15,2).Value2 = $desc
$ws.Cells.Item(15,3).Value2 = $grade
$ws.Cells.Item(15,4).Value2 = 264.5
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value2 = "01.08.2025"
$ws.Cells.Item(15,6).Value2 = ""

$ws.Cells.Item(16,1).Value2 = 15
$ws.Cells.Item(16,2).Value2 = $desc
$ws.Cells.Item(16,3).Value2 = $grade
$ws.Cells.Item(16,4).Value2 = 266.25
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value2 = "29.07.2025"
$ws.Cells.Item(16,6).Value2 = ""

$ws.Cells.Item(17,1).Value2 = 14
$ws.Cells.Item(17,2).Value2 = $desc
$ws.Cells.Item(17,3).Value2 = $grade
$ws.Cells.Item(17,4).Value2 = 268.5
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value2 = "26.07.2025"
$ws.Cells.Item(17,6).Value2 = ""

$ws.Cells.Item(18,1).Value2 = 13
$ws.Cells.Item(18,2).Value2 = $desc
$ws.Cells.Item(18,3).Value2 = $grade
$ws.Cells.Item(18,4).Value2 = 267
$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value2 = "22.07.2025"
$ws.Cells.Item(18,6).Value2 = ""

$ws.Cells.Item(19,1).Value2 = 12
$ws.Cells.Item(19,2).Value2 = $desc
$ws.Cells.Item(19,3).Value2 = $grade
$ws.Cells.Item(19,4).Value2 = 261.5
$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value2 = "19.07.2025"
$ws.Cells.Item(19,6).Value2 = ""

$ws.Cells.Item(20,1).Value2 = 11
$ws.Cells.Item(20,2).Value2 = $desc
$ws.Cells.Item(20,3).Value2 = $grade
$ws.Cells.Item(20,4).Value2 = 258
$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value2 = "17.07.2025"
$ws.Cells.Item(20,6).Value2 = ""

$ws.Cells.Item(21,1).Value2 = 10
$ws.Cells.Item(21,2).Value2 = $desc
$ws.Cells.Item(21,3).Value2 = $grade
$ws.Cells.Item(21,4).Value2 = 261.25
$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,5).Value2 = "11.07.2025"
$ws.Cells.Item(21,6).Value2 = ""

$ws.Cells.Item(22,1).Value2 = 9
$ws.Cells.Item(22,2).Value2 = $desc
$ws.Cells.Item(22,3).Value2 = $grade
$ws.Cells.Item(22,4).Value2 = 258.5
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value2 = "05.07.2025"
$ws.Cells.Item(22,6).Value2 = ""

$ws.Cells.Item(23,1).Value2 = 8
$ws.Cells.Item(23,2).Value2 = $desc
$ws.Cells.Item(23,3).Value2 = $grade
$ws.Cells.Item(23,4).Value2 = 260.75
$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value2 = "02.07.2025"
$ws.Cells.Item(23,6).Value2 = ""

$ws.Cells.Item(24,1).Value2 = 7
$ws.Cells.Item(24,2).Value2 = $desc
$ws.Cells.Item(24,3).Value2 = $grade
$ws.Cells.Item(24,4).Value2 = 263.25
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value2 = "28.06.2025"
$ws.Cells.Item(24,6).Value2 = ""

$ws.Cells.Item(25,1).Value2 = 6
$ws.Cells.Item(25,2).Value2 = $desc
$ws.Cells.Item(25,3).Value2 = $grade
$ws.Cells.Item(25,4).Value2 = 261.75
$ws.Cells.Item(25,5).NumberFormat = "@"
$ws.Cells.Item(25,5).Value2 = "26.06.2025"
$ws.Cells.Item(25,6).Value2 = ""

$ws.Cells.Item(26,1).Value2 = 5
$ws.Cells.Item(26,2).Value2 = $desc
$ws.Cells.Item(26,3).Value2 = $grade
$ws.Cells.Item(26,4).Value2 = 264
$ws.Cells.Item(26,5).NumberFormat = "@"
$ws.Cells.Item(26,5).Value2 = "25.06.2025"
$ws.Cells.Item(26,6).Value2 = ""

$ws.Cells.Item(27,1).Value2 = 4
$ws.Cells.Item(27,2).Value2 = $desc
$ws.Cells.Item(27,3).Value2 = $grade
$ws.Cells.Item(27,4).Value2 = 268.75
$ws.Cells.Item(27,5).NumberFormat = "@"
$ws.Cells.Item(27,5).Value2 = "24.06.2025"
$ws.Cells.Item(27,6).Value2 = ""

$ws.Cells.Item(28,1).Value2 = 3
$ws.Cells.Item(28,2).Value2 = $desc
$ws.Cells.Item(28,3).Value2 = $grade
$ws.Cells.Item(28,4).Value2 = 262.25
$ws.Cells.Item(28,5).NumberFormat = "@"
$ws.Cells.Item(28,5).Value2 = "19.06.2025"
$ws.Cells.Item(28,6).Value2 = ""

$ws.Cells.Item(29,1).Value2 = 2
$ws.Cells.Item(29,2).Value2 = $desc
$ws.Cells.Item(29,3).Value2 = $grade
$ws.Cells.Item(29,4).Value2 = 260
$ws.Cells.Item(29,5).NumberFormat = "@"
$ws.Cells.Item(29,5).Value2 = "18.06.2025"
$ws.Cells.Item(29,6).Value2 = ""

$ws.Cells.Item(30,1).Value2 = 1
$ws.Cells.Item(30,2).Value2 = $desc
$ws.Cells.Item(30,3).Value2 = $grade
$ws.Cells.Item(30,4).Value2 = 256.5
$ws.Cells.Item(30,5).NumberFormat = "@"
$ws.Cells.Item(30,5).Value2 = "12.06.2025"
$ws.Cells.Item(30,6).Value2 = ""

# --- Re-add hyperlinks for F2:F10 ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-02-september-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-27-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-26-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-19-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-14-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf") | Out-Null

# --- Restore original cell formatting (style 3 / style 4) on every cell we touched ---
# (Value2 assignment itself preserves pre-existing style, but brand-new row 30 cells start
#  with no style, and NumberFormat="@"/Hyperlinks.Add override style on E/F columns, so
#  re-apply the captured style templates uniformly to every data cell, A:F, rows 2-30.)
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A2:C2").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(2,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E2:F2").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A3:C3").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(3,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E3:F3").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A4:C4").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(4,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E4:F4").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A5:C5").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(5,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E5:F5").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A6:C6").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(6,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E6:F6").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A7:C7").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(7,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E7:F7").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A8:C8").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(8,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E8:F8").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A9:C9").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(9,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E9:F9").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A10:C10").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(10,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E10:F10").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A11:C11").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(11,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E11:F11").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A12:C12").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(12,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E12:F12").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A13:C13").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(13,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E13:F13").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A14:C14").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(14,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E14:F14").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A15:C15").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(15,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E15:F15").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A16:C16").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(16,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E16:F16").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A17:C17").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(17,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E17:F17").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A18:C18").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(18,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E18:F18").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A19:C19").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(19,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E19:F19").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A20:C20").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(20,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E20:F20").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A21:C21").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(21,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E21:F21").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A22:C22").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(22,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E22:F22").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A23:C23").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(23,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E23:F23").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A24:C24").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(24,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E24:F24").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A25:C25").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(25,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E25:F25").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A26:C26").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(26,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E26:F26").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A27:C27").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(27,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E27:F27").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A28:C28").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(28,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E28:F28").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A29:C29").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(29,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E29:F29").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("A30:C30").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,11).Copy() | Out-Null
$ws.Cells.Item(30,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,10).Copy() | Out-Null
$ws.Range("E30:F30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Remove scratch template cells ---
$ws.Cells.Item(1,10).Clear() | Out-Null
$ws.Cells.Item(1,11).Clear() | Out-Null

# --- Column F width: stored XML width 79 -> 80 (runtime adds +5/6 when ColumnWidth is set explicitly) ---
$ws.Columns.Item(6).ColumnWidth = 80 - (5/6)

$wb.Save()